$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "307.88") stay text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.524.18"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.638.98"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "307.88"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.3770"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "52.49"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "0.3653"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "1.271"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "0.08184"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "23.00"
$ws.Range("E13").Value = "  +1.92%  "
$ws.Range("D14").Value = "6.641"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "0.00001279"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "7.401"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "1.637.29"
$ws.Range("D18").Value = "94.71"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "0.06938"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "6.563"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "23.522.97"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "12.82"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "3.083"
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("D26").Value = "2.424"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "21.29"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "151.39"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "5.357"
$ws.Range("D30").Value = "135.59"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").Value = "2.379"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "1.824.71"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "0.9728"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "0.02825"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("D36").Value = "10.33"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "0.07361"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").Value = "0.2549"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").Value = "6.183"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Value = "0.08890"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").Value = "1.382"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "0.7111"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("D43").Value = "12.51"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "16.23"
$ws.Range("E44").Value = "  +5.39%  "
$ws.Range("D45").Value = "0.6547"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "2.343"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "4.042"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "0.07977"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "129.60"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "1.211"
$ws.Range("E51").Value = "  +0.64%  "
